$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TLRD")

# Insert a new column before D; existing D:K data shifts to E:L
$ws.Columns("D").Insert(1)

# New column D should look like the (shifted) former column D, now at E:
# same cell formatting / number formats / column width.
$ws.Columns("E").Copy()
$ws.Columns("D").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Columns("D").ColumnWidth = $ws.Columns("E").ColumnWidth

# Populate new column D with the newest quarter (period ending 2018-11-03) figures
$ws.Range("D7").Value = 43407
$ws.Range("D8").Value = 812700
$ws.Range("D9").Value = 450000
$ws.Range("D10").Value = 362700
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 33400
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 767100
$ws.Range("D18").Value = 45600
$ws.Range("D20").Value = 300
$ws.Range("D21").Value = 71300
$ws.Range("D22").Value = 18800
$ws.Range("D23").Value = 27100
$ws.Range("D24").Value = 13300
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 13900
$ws.Range("D27").Value = 13900
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -300
$ws.Range("D33").Value = 13900
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 13900
$ws.Range("D38").Value = 43407
$ws.Range("D41").Value = 68400
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 81700
$ws.Range("D44").Value = 875000
$ws.Range("D45").Value = 69500
$ws.Range("D46").Value = 1094600
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 430900
$ws.Range("D49").Value = 244300
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 119800
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 1889500
$ws.Range("D57").Value = 236000
$ws.Range("D58").Value = 9000
$ws.Range("D59").Value = 318600
$ws.Range("D60").Value = 563600
$ws.Range("D61").Value = 1167900
$ws.Range("D62").Value = 148600
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 1880100
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = -465000
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 9500
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43407
$ws.Range("D81").Value = 13900
$ws.Range("D83").Value = 25400
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 79800
$ws.Range("D91").Value = -22300
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -22300
$ws.Range("D96").Value = -9100
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -57400
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = 200
